$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Authors")

# New row of data (row 3 is intentionally left empty / not written, matching
# the "empty rows are not generated in the XML" fix - the row after the
# last written row, 3, is skipped and the new data lands on row 4).
$ws1.Range("A4").Value = "Small Bird"
$ws1.Range("B4").Value = "Cleaner"

# Move the selection down to B5 and make "Authors" the active/selected sheet.
$ws1.Range("B5").Select()
$ws1.Activate()
